# Restored from revision of admin on 11/04/2020 05:34:46 PM.TEST
# Author: admin. Type: SAVE.
#
# Rule row 10 (R30): the "Integer min" condition value (column C / C10) was
# changed from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C10").Value = 1
